$wb = $excel.ActiveWorkbook
$xlPasteFormats = -4122

# --- Step 1: duplicate the "2022-Q3" sheet, placing the copy right after it.
# The duplicate keeps the old (Q3) figures and becomes the archived "2022-Q3"
# sheet, while the original sheet is updated in place to become "2022-Q4".
$srcSheet = $wb.Worksheets.Item("2022-Q3")
$q3Index = $srcSheet.Index
$srcSheet.Copy($null, $srcSheet)
$newSheet = $wb.Worksheets.Item($q3Index + 1)
$newSheet.Name = "2022-Q3 (archived)"

# --- Step 2: rename the original sheet to "2022-Q4" and refresh its figures.
$q4Sheet = $wb.Worksheets.Item("2022-Q3")
$q4Sheet.Name = "2022-Q4"

$q4Sheet.Range("D2:G2").NumberFormat = "@"
$q4Sheet.Range("D2").Value = "0.23"
$q4Sheet.Range("E2").Value = "94.47"
$q4Sheet.Range("F2").Value = "2.61"
$q4Sheet.Range("G2").Value = "0.0060"
$q4Sheet.Range("H2").Value = 3

# --- Step 3: rename the duplicate back to "2022-Q3" (keeps the old figures).
$newSheet.Name = "2022-Q3"

# --- Step 4: update the "总计" (summary) sheet with the new quarter's row,
# pushing the existing "2021-Q3" total down one row.
$sumSheet = $wb.Worksheets.Item("总计")
$sumSheet.Range("B2").Value = "2022-Q4"

$sumSheet.Range("A4").Value = 2
$sumSheet.Range("A3").Copy()
$sumSheet.Range("A4").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

$sumSheet.Range("B4").Value = "2021-Q3"
$sumSheet.Range("C4").Value = 2
$sumSheet.Range("D4").Value = 0.18

$sumSheet.Range("B3").Value = "2022-Q3"
$sumSheet.Range("C3").Value = 1
$sumSheet.Range("D3").Value = 0.01

# --- Step 5: restore the originally-active "2021-Q3" tab (the sheet Copy in
# step 1 shifts the active tab onto the newly inserted sheet).
$wb.Worksheets.Item("2021-Q3").Activate()
